{"js": "// The document has paragraphs whose <w:pPr> contains a stray, leading\n// <w:pStyle w:val=\"Compact\"/> in front of the \"real\" paragraph style\n// (Word/OOXML lets the last <w:pStyle> win, so these paragraphs already\n// *display* with their intended style, e.g. \"Casebook Title\"). This\n// script removes that stray \"Compact\" style everywhere it shadows a more\n// specific style, and for the one paragraph whose *only* style is\n// \"Compact\" (a footnote-text paragraph), switches it to \"Body Text\" -\n// matching the target diff exactly.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/style\");\nawait context.sync();\n\n// Styles that (in this document) are always layered on top of a stray\n// leading \"Compact\" pStyle. Re-assigning the paragraph's own displayed\n// style forces Word to collapse the duplicate pStyle entries down to a\n// single, correct one.\nconst collapseStyles = new Set([\n  \"Casebook Title\",\n  \"Casebook Subtitle\",\n  \"Section Number\",\n  \"Section Title\",\n  \"Section Subtitle\",\n  \"Resource Number\",\n  \"Resource Title\",\n  \"Resource Link\",\n]);\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  const style = paragraph.style;\n\n  if (style === \"Compact\") {\n    // The lone \"my note\" paragraph: Compact -> Body Text.\n    paragraph.style = \"Body Text\";\n  } else if (collapseStyles.has(style)) {\n    // Re-apply the same (already-effective) style so Word drops the\n    // now-redundant leading \"Compact\" pStyle.\n    paragraph.style = style;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document has paragraphs whose paragraph properties contain a\n# stray, leading Compact paragraph style in front of the \"real\" style\n# (OOXML lets the last <w:pStyle> win, so these paragraphs already\n# *display*/report their intended style, e.g. \"Casebook Title\"). This\n# script removes that stray \"Compact\" style everywhere it shadows a\n# more specific style, and for the one paragraph whose *only* style is\n# \"Compact\" (a footnote-text paragraph), switches it to \"Body Text\" -\n# matching the target diff exactly.\n\n$d = $word.ActiveDocument\n\n# Styles that (in this document) are always layered on top of a stray\n# leading \"Compact\" style. Re-assigning the paragraph's own reported\n# style forces Word to collapse the duplicate style entries down to a\n# single, correct one.\n$collapseStyles = @(\n    \"Casebook Title\",\n    \"Casebook Subtitle\",\n    \"Section Number\",\n    \"Section Title\",\n    \"Section Subtitle\",\n    \"Resource Number\",\n    \"Resource Title\",\n    \"Resource Link\"\n)\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs($i)\n    $styleName = $p.Style.NameLocal\n\n    if ($styleName -eq \"Compact\") {\n        # The lone \"my note\" paragraph: Compact -> Body Text.\n        $p.Style = \"Body Text\"\n    } elseif ($collapseStyles -contains $styleName) {\n        # Re-apply the same (already-effective) style so Word drops the\n        # now-redundant leading \"Compact\" style.\n        $p.Style = $styleName\n    }\n}\n"}
